$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('463', 'Bag Wax - 6#', '1', '$62.72', '$62.72'),
    @('18356', 'Container - Plastic Clamshell (6")', '1', '$58.23', '$58.23'),
    @('34365', 'Towel Disp Rolls (IB)', '2', '$53.88', '$107.76'),
    @('26215', 'Glove Latex - Med', '1', '$38.62', '$38.62'),
    @('26216', 'Glove Latex - Large', '1', '$38.62', '$38.62'),
    @('21094', 'Lid - Soup (6/16 oz)', '1', '$39.58', '$39.58'),
    @('27455', 'Container - Soup (16oz)', '3', '$45.62', '$136.86'),
    @('38018', 'Cont Salad Bowl - 24oz (paper)', '1', '$41.92', '$41.92'),
    @('14909', 'Lid - Deli (6/32oz)', '2', '$24.82', '$49.64'),
    @('19608', 'Cup - Cold (24oz)', '1', '$70.67', '$70.67'),
    @('1847', 'Wrap Sub 18x18 (Cold)', '1', '$102.85', '$102.85'),
    @('39808', 'Wrap Cushion 14x16 (Hot)', '2', '$38.65', '$77.30'),
    @('26214', 'Glove Latex - Small', '1', '$38.62', '$38.62')
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $cell = $ws.Cells.Item($rowNum, $c + 1)
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$c]
        $cell.Style = "Normal"
    }
}
